$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "HMS All Doctor Display" API row (row 10) -----------------
# (values are written in the same order the shared-strings table ends up
# with: url, name, then output)
$ws.Range("C10").Value = "http://127.0.0.1:8000/staff/displayDoctor/"
$ws.Range("A10").Value = "HMS All Doctor Display"
$ws.Range("B10").Value = "POST"
$ws.Range("E10").Value = "[{""doctorid"": 1, ""name"": ""ASIM THAHA AZEEZ"", ""username"": ""asimthaha"", ""speciality"": ""Nephrologist"", ""startYear"": 2000,`n""qualification"": ""MBBS, MD"", ""role"": ""Doctor"", ""password"": ""1234""}]"

# Copy the formatting from row 9 (same layout: A/B bold-centered, C hyperlink
# style, D blank, E centered-wrap) onto the new row so it matches the rest
# of the table.
$ws.Range("A9:C9").Copy() | Out-Null
$ws.Range("A10:C10").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("E9").Copy() | Out-Null
$ws.Range("E10").PasteSpecial(-4122) | Out-Null        # xlPasteFormats
$excel.CutCopyMode = 0

# Turn the API url into a real hyperlink, then restore the formatting once
# more so the cell keeps the table's standard hyperlink look.
$ws.Hyperlinks.Add($ws.Range("C10"), "http://127.0.0.1:8000/staff/displayDoctor/") | Out-Null
$ws.Range("C9").Copy() | Out-Null
$ws.Range("C10").PasteSpecial(-4122) | Out-Null         # xlPasteFormats
$excel.CutCopyMode = 0

# --- Row heights (re-fit now that the table has an extra row) --------------
$ws.Rows(2).RowHeight = 57.6
$ws.Rows(3).RowHeight = 28.8
$ws.Rows(4).RowHeight = 57.6
$ws.Rows(5).RowHeight = 43.2
$ws.Rows(6).AutoFit()   # content now fits the standard row height
$ws.Rows(7).AutoFit()   # content now fits the standard row height
$ws.Rows(8).RowHeight = 144
$ws.Rows(9).RowHeight = 72
$ws.Rows(10).RowHeight = 72

# --- View state: scroll back to the top and select G5 ----------------------
$ws.Application.Goto($ws.Range("A1"))
$ws.Range("G5").Select()
